$wb = $excel.ActiveWorkbook

# --- Rename sheet tabs ---
$wb.Worksheets.Item(1).Name = "GNG_TO-1651255486796629"
$wb.Worksheets.Item(2).Name = "NB_TO-16512554887478058"
$wb.Worksheets.Item(3).Name = "RS_TO-16512554887498028"
$wb.Worksheets.Item(4).Name = "TOL_TO-16512554888117676"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16512554888900013"

# --- Sheet1 (GNG) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16512554867546308.csv"
$ws1.Range("B3").Value = "GNG_stims-1651255486778631.csv"
$ws1.Range("B4").Value = "go_stims-16512554867796311.csv"
$ws1.Range("B5").Value = "GNG_stims-1651255486795629.csv"

# --- Sheet2 (NB) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "ZB-match_0-16512554870566304.csv"
$ws2.Range("B3").Value = "OB-1651255487461635.csv"
$ws2.Range("B4").Value = "ZB-match_4-16512554868556652.csv"
$ws2.Range("B5").Value = "TB-16512554878496318.csv"
$ws2.Range("B6").Value = "TB-1651255488212629.csv"
$ws2.Range("B7").Value = "ZB-match_4-16512554870216627.csv"
$ws2.Range("B8").Value = "OB-16512554871066618.csv"
$ws2.Range("B9").Value = "OB-1651255487082663.csv"
$ws2.Range("B10").Value = "TB-16512554887287996.csv"

# --- Sheet4 (TOL) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16512554887638009.csv"
$ws4.Range("B3").Value = "ZM_stims-16512554887517667.csv"
$ws4.Range("B4").Value = "MM_stims-16512554887947762.csv"
$ws4.Range("B5").Value = "ZM_stims-16512554887647672.csv"
$ws4.Range("B6").Value = "MM_stims-16512554888107657.csv"
$ws4.Range("B7").Value = "ZM_stims-16512554887957673.csv"

# --- Sheet5 (vSAT) ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16512554888579721.csv"
$ws5.Range("B3").Value = "SAT_stims-1651255488842002.csv"
$ws5.Range("B4").Value = "vSAT_stims-16512554888740048.csv"
$ws5.Range("B5").Value = "SAT_stims-16512554888157709.csv"
